# Applies the "smoy_winners" update:
#   - renames header columns year_x -> season_ending_year_x, year_y -> season_ending_year_y
#   - back-fills the birth_year column (Q) for every player row
#   - appends a new calendar_year column (AY) derived from the season_y column (AV)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the two "year" headers to their clearer names -------------------
$ws.Range("A1").Value = "season_ending_year_x"
$ws.Range("O1").Value = "season_ending_year_y"

# --- Add the new calendar_year header (col 51 / AY), matching the look of
#     the other header cells (bold, centered, thin border) ------------------
$ws.Cells.Item(1, 51).Value = "calendar_year"
$ws.Cells.Item(1, 51).Font.Bold = $true
$ws.Cells.Item(1, 51).HorizontalAlignment = -4108
$ws.Cells.Item(1, 51).VerticalAlignment = -4160
$ws.Cells.Item(1, 51).Borders.LineStyle = 1

# --- Fill in birth_year (col 17 / Q) and calendar_year (col 51 / AY) for
#     every data row, derived from season_y (col 48 / AV, e.g. "2023-24")
#     and age_y (col 19 / S) -------------------------------------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 48).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $seasonY = $ws.Cells.Item($row, 48).Value()
    $ageY = $ws.Cells.Item($row, 19).Value()

    $startYear = [int]($seasonY.Substring(0, 4))
    $calendarYear = $startYear + 1
    $birthYear = $calendarYear - [int]$ageY

    $ws.Cells.Item($row, 17).Value = $birthYear
    $ws.Cells.Item($row, 51).Value = $calendarYear
}
